$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = -7.575
$ws.Range("D10").Value = -7.999
$ws.Range("D12").Value = -7.886000000000001
$ws.Range("E13").Value = 12.583
$ws.Range("D18").Value = -8.167
